# Add 7 new rows (331-337) of landscaping observations for 6/26/2025
# to the bottom of the data table on Sheet1, mirroring the layout of
# the existing rows (A:T) and extending the shared ABS() formula in
# column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 331
$lastNewRow  = 337

# --- Column A: date serial for every new row, with the same date
#     number-format as the existing rows (copy format from A330). ---
$ws.Range("A330").Copy() | Out-Null
$ws.Range("A$firstNewRow`:A$lastNewRow").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("A$firstNewRow`:A$lastNewRow").Value2 = 45834

# --- Columns B through T for the 7 new rows. ---
$data = @(
    @("Flowering",    "Large",  72, 91, 0.22, 0.2,  "Yes", 2, "Bright",  9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0),
    @("Nonflowering", "Medium", 72, 91, 0.22, 0.3,  "Yes", 3, "Bright",  9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0),
    @("Nonflowering", "Small",  72, 91, 0.22, 0.25, "Yes", 3, "Dark",    9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0),
    @("Nonflowering", "Medium", 72, 91, 0.22, 0.4,  "Yes", 3, "Bright",  9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0),
    @("Nonflowering", "Medium", 72, 91, 0.22, 0.5,  "Yes", 3, "Dark",    9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0),
    @("Nonflowering", "Large",  72, 91, 0.22, 0.8,  "Yes", 4, "Neutral", 9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0),
    @("Tree",          "Medium",72, 91, 0.22, 1.2,  "Yes", 1, "Bright",  9, 0.59, 74, 30.01, 8, 0.5, 9.9, 51, 0)
)

$rowCount = $data.Count
$colCount = $data[0].Count   # B..T minus F = 18 values (B,C,D,E,G,H,I,J,K,L,M,N,O,P,Q,R,S,T)

$arr = New-Object 'object[,]' $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

# B:E (4 cols)
$bToE = New-Object 'object[,]' $rowCount, 4
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt 4; $c++) { $bToE[$r, $c] = $arr[$r, $c] }
}
$ws.Range("B$firstNewRow`:E$lastNewRow").Value2 = $bToE

# G:T (14 cols) -- everything after the D/E/F group
$gToT = New-Object 'object[,]' $rowCount, 14
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 4; $c -lt $colCount; $c++) { $gToT[$r, $c - 4] = $arr[$r, $c] }
}
$ws.Range("G$firstNewRow`:T$lastNewRow").Value2 = $gToT

# --- Column F: extend the ABS(D-E) formula down through the new rows. ---
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Range("F$r").Formula = "=ABS(D$r-E$r)"
}

# --- Move the selection to where the user ended up after the paste
#     (I338, right after the new block). ---
$ws.Range("I338").Select() | Out-Null

$wb.Application.CalculateFullRebuild()
